$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E9").Value = 17.35510000000001
$ws.Range("E18").Value = 17.88670000000002
$ws.Range("E20").Value = 16.05899999999999
$ws.Range("E27").Value = 16.70029999999999
$ws.Range("E69").Value = 17.31520000000002
$ws.Range("E76").Value = 16.14359999999999
$ws.Range("E82").Value = 16.6852
